# Case_2_35/line/parallel.xlsx — "contingencies with rene fine"
# Adds two new columns (P, Q) to the existing A1:O25 table and flips the
# values in columns I, K, M, O for every data row (2-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- extend the header row (row 1) with two new columns -------------------
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Match the header formatting used by the rest of row 1 (bold, centered,
# thin border all around) so the new header cells pick up the same style
# as B1:O1.
$hdr = $ws.Range("P1:Q1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# --- update data rows 2-25 --------------------------------------------------
for ($r = 2; $r -le 25; $r++) {
    # Swap the 1/2 values in columns I, K, M, O.
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1

    # New columns P and Q both carry the value 2.
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
